$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold plain-looking decimal numbers (e.g. "213.01") that
# Excel's type-inference would otherwise silently convert to a Number.
# The source data models Price as text (inline strings), so force the
# number format to Text first to keep them as strings after assignment.
$textCells = @("D5","D10","D11","D17","D19","D22","D23","D25","D28","D30","D34","D37","D38","D40","D41","D45","D46","D47","D48","D49")
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.653.70"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.24"
$ws.Range("E3").Value = "  +0.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.01"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +3.23%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.13%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +2.13%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.22%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  +2.02%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +3.53%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "1.861.16"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13 - Wrapped Ether
$ws.Range("D13").Value = "1.625.09"
$ws.Range("E13").Value = "  +0.26%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.97%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.25%  "

# Row 16 - Wrapped BTC
$ws.Range("D16").Value = "26.665.43"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.32"
$ws.Range("E17").Value = "  +1.64%  "

# Row 18 - Shiba Inu
$ws.Range("E18").Value = "  +1.79%  "

# Row 19 - Bitcoin Cash
$ws.Range("D19").Value = "218.81"
$ws.Range("E19").Value = "  +8.44%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.08%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.09%  "

# Row 22 & 23 - Chainlink / Avalanche swap places
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "9.43"
$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.19"
$ws.Range("E23").Value = "  +2.47%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +2.71%  "

# Row 25 - Monero
$ws.Range("D25").Value = "148.26"
$ws.Range("E25").Value = "  +2.35%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.09%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.35%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "6.89"
$ws.Range("E28").Value = "  +5.04%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +1.96%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -2.40%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.41%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +3.85%  "

# Row 33 - InternetComputer (DFINITY)
$ws.Range("E33").Value = "  +1.65%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  +0.97%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.08%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.212.01"
$ws.Range("E36").Value = "  +2.69%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  +5.41%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.812"
$ws.Range("E38").Value = "  +0.38%  "

# Row 39 - PaxDollar
$ws.Range("E39").Value = "  +0.14%  "

# Row 40 - ImmutableX
$ws.Range("D40").Value = "0.502"
$ws.Range("E40").Value = "  +1.26%  "

# Row 41 - MXToken
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  -1.42%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  +1.40%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  +0.79%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.772.92"
$ws.Range("E44").Value = "  +0.87%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.77"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "1.56"
$ws.Range("E46").Value = "  +1.00%  "

# Row 47 - Aave
$ws.Range("D47").Value = "54.71"
$ws.Range("E47").Value = "  +1.68%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +1.70%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "7.64"
$ws.Range("E49").Value = "  +4.82%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +0.30%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.22%  "
